$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, shifting rows 111:131 down to 112:132
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with data
$ws.Range("A111").Value = 8
$ws.Range("B111").Value = "Terminal La Palmera de La Serena"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = 44644
$ws.Range("E111").Value = 4
$ws.Range("F111").Value = 100112044
$ws.Range("G111").Value = "Perejil"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 2400
$ws.Range("K111").Value = 2300
$ws.Range("L111").Value = 2500
$ws.Range("M111").Value = 2400
$ws.Range("N111").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O111").Value = "Provincia del Elquí"
$ws.Range("P111").Value = 1600
$ws.Range("Q111").Value = 1.5
$ws.Range("R111").Value = "Hortaliza"
